$d = $word.ActiveDocument

# --- 1. Locate the "Code: Siehe Java in Github" paragraph ---------------
$codeParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Code: Siehe Java in Github*") {
        $codeParaIndex = $i
        break
    }
}

# --- 2. Remove the old (end-of-document) _GoBack bookmark ---------------
try {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
} catch {
}

# --- 3. Insert the hyperlink right before that paragraph's text ---------
$codePara = $d.Paragraphs.Item($codeParaIndex)
$insertionPoint = $codePara.Range
$insertionPoint.Collapse(1)

$url = "https://lucid.app/lucidchart/3eda4911-0192-4dd9-a8c6-ceef8a33c2e6/edit?viewport_loc=-1892%2C-591%2C2558%2C1130%2C0_0&invitationId=inv_3916ee22-1b50-4787-ac64-9681b245df12"
$hyperlink = $d.Hyperlinks.Add($insertionPoint, $url)
$hyperlink.Range.Font.Size = 12

# --- 4. Split the paragraph so the hyperlink sits in its own paragraph --
$splitPoint = $d.Range($hyperlink.Range.End, $hyperlink.Range.End)
$splitPoint.InsertParagraphAfter()

# --- 5. Re-create the _GoBack bookmark right before "Code: ..." text ----
$codePara = $d.Paragraphs.Item($codeParaIndex + 1)
$bookmarkRange = $codePara.Range
$bookmarkRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- 6. Add the "Hyperlink" character style ------------------------------
$hlStyle = $d.Styles.Add("Hyperlink", 2)
$hlStyle.BaseStyle = $d.Styles.Item("Absatz-Standardschriftart")
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.Font.Underline = 1
$hlStyle.Font.Color = 0xC16305

Write-Host "done"
